$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cell H1, matching the style used by the other header cells (bold, bordered, centered)
$ws.Range("H1").Value = "Save"
$ws.Range("H1").Font.Bold = $true
$ws.Range("H1").HorizontalAlignment = -4108  # xlCenter
$ws.Range("H1").VerticalAlignment = -4160    # xlTop
$ws.Range("H1").Borders.LineStyle = 1        # xlContinuous
$ws.Range("H1").Borders.Weight = 2           # xlThin

# New data cells
$ws.Range("H2").Value = 1
$ws.Range("H3").Value = 0
